$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 83, shifting the existing rows 83-88 down to 84-89.
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new weekly record.
$ws.Cells.Item(83, 1).Value = 5
$ws.Cells.Item(83, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(83, 3).Value = "Maule"
$ws.Cells.Item(83, 4).Value = 44585
$ws.Cells.Item(83, 5).Value = 7
$ws.Cells.Item(83, 6).Value = 100112001
$ws.Cells.Item(83, 7).Value = "Berenjena"
$ws.Cells.Item(83, 8).Value = "Sin especificar"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 200
$ws.Cells.Item(83, 11).Value = 7000
$ws.Cells.Item(83, 12).Value = 7000
$ws.Cells.Item(83, 13).Value = 7000
$ws.Cells.Item(83, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(83, 15).Value = "Región del Maule"
$ws.Cells.Item(83, 16).Value = 117
$ws.Cells.Item(83, 17).Value = 60
$ws.Cells.Item(83, 18).Value = "Hortaliza"
